# Updates the TPM-derived NATMI edge metrics for the Ctf1-Il6st LR pair sheet.
# Re-labels the "Target cluster" column (D) using the refreshed shared-string
# catalogue and rewrites the recomputed expression / specificity values for rows 2-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.2720863333333334
$ws.Range("H2").Value = 0.8162590000000001
$ws.Range("I2").Value = 0.1601964374275199
$ws.Range("J2").Value = 0.2007414748833069
$ws.Range("M2").Value = 55.783591
$ws.Range("N2").Value = 167.350773
$ws.Range("O2").Value = 0.2332214199005771
$ws.Range("P2").Value = 0.2394371967339281
$ws.Range("Q2").Value = 15.17795273535634
$ws.Range("R2").Value = 136.601574618207
$ws.Range("S2").Value = 0.03736124059986015
$ws.Range("T2").Value = 0.04806497601429324

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.2720863333333334
$ws.Range("H3").Value = 0.8162590000000001
$ws.Range("I3").Value = 0.1601964374275199
$ws.Range("J3").Value = 0.2007414748833069
$ws.Range("O3").Value = 0.5297503589663128
$ws.Range("P3").Value = 0.5438691736537713
$ws.Range("Q3").Value = 34.4759324137399
$ws.Range("R3").Value = 310.283391723659
$ws.Range("S3").Value = 0.08486412023235314
$ws.Range("T3").Value = 0.1091771000628234

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.2720863333333334
$ws.Range("H4").Value = 0.8162590000000001
$ws.Range("I4").Value = 0.1601964374275199
$ws.Range("J4").Value = 0.2007414748833069
$ws.Range("M4").Value = 23.03749833333333
$ws.Range("N4").Value = 69.112495
$ws.Range("O4").Value = 0.09631574403765399
$ws.Range("P4").Value = 0.09888273454277752
$ws.Range("Q4").Value = 6.268188450689445
$ws.Range("R4").Value = 56.413696056205
$ws.Range("S4").Value = 0.01542943906301306
$ws.Range("T4").Value = 0.01984986597261168

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.2720863333333334
$ws.Range("H5").Value = 0.8162590000000001
$ws.Range("I5").Value = 0.1601964374275199
$ws.Range("J5").Value = 0.2007414748833069
$ws.Range("M5").Value = 18.627865
$ws.Range("N5").Value = 37.25573
$ws.Range("O5").Value = 0.07787983970082285
$ws.Range("P5").Value = 0.05330365312071852
$ws.Range("Q5").Value = 5.068387485678334
$ws.Range("R5").Value = 30.41032491407
$ws.Range("S5").Value = 0.01247607286749815
$ws.Range("T5").Value = 0.01070025394412122

# Row 6
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.2720863333333334
$ws.Range("H6").Value = 0.8162590000000001
$ws.Range("I6").Value = 0.1601964374275199
$ws.Range("J6").Value = 0.2007414748833069
$ws.Range("M6").Value = 15.028766
$ws.Range("N6").Value = 45.086298
$ws.Range("O6").Value = 0.06283263739463307
$ws.Range("P6").Value = 0.06450724194880479
$ws.Range("Q6").Value = 4.089121835464667
$ws.Range("R6").Value = 36.802096519182
$ws.Range("S6").Value = 0.01006556466479538
$ws.Range("T6").Value = 0.0129492788894574

# Row 7
$ws.Range("D7").Value = "ECs"
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.397225
$ws.Range("H7").Value = 1.191675
$ws.Range("I7").Value = 0.2338744069853316
$ws.Range("J7").Value = 0.2930670253945926
$ws.Range("M7").Value = 55.783591
$ws.Range("N7").Value = 167.350773
$ws.Range("O7").Value = 0.2332214199005771
$ws.Range("P7").Value = 0.2394371967339281
$ws.Range("Q7").Value = 22.158636934975
$ws.Range("R7").Value = 199.427732414775
$ws.Range("S7").Value = 0.05454452127552448
$ws.Range("T7").Value = 0.07017114701563215

# Row 8
$ws.Range("D8").Value = "FAPs"
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.397225
$ws.Range("H8").Value = 1.191675
$ws.Range("I8").Value = 0.2338744069853316
$ws.Range("J8").Value = 0.2930670253945926
$ws.Range("O8").Value = 0.5297503589663128
$ws.Range("P8").Value = 0.5438691736537713
$ws.Range("Q8").Value = 50.33219451074167
$ws.Range("R8").Value = 452.989750596675
$ws.Range("S8").Value = 0.123895051053513
$ws.Range("T8").Value = 0.1593901209265259

# Row 9
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.397225
$ws.Range("H9").Value = 1.191675
$ws.Range("I9").Value = 0.2338744069853316
$ws.Range("J9").Value = 0.2930670253945926
$ws.Range("M9").Value = 23.03749833333333
$ws.Range("N9").Value = 69.112495
$ws.Range("O9").Value = 0.09631574403765399
$ws.Range("P9").Value = 0.09888273454277752
$ws.Range("Q9").Value = 9.151070275458332
$ws.Range("R9").Value = 82.359632479125
$ws.Range("S9").Value = 0.02252578752015731
$ws.Range("T9").Value = 0.02897926887533494

# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.397225
$ws.Range("H10").Value = 1.191675
$ws.Range("I10").Value = 0.2338744069853316
$ws.Range("J10").Value = 0.2930670253945926
$ws.Range("M10").Value = 18.627865
$ws.Range("N10").Value = 37.25573
$ws.Range("O10").Value = 0.07787983970082285
$ws.Range("P10").Value = 0.05330365312071852
$ws.Range("Q10").Value = 7.399453674625
$ws.Range("R10").Value = 44.39672204775
$ws.Range("S10").Value = 0.01821410132614263
$ws.Range("T10").Value = 0.01562154306275417

# Row 11
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.397225
$ws.Range("H11").Value = 1.191675
$ws.Range("I11").Value = 0.2338744069853316
$ws.Range("J11").Value = 0.2930670253945926
$ws.Range("M11").Value = 15.028766
$ws.Range("N11").Value = 45.086298
$ws.Range("O11").Value = 0.06283263739463307
$ws.Range("P11").Value = 0.06450724194880479
$ws.Range("Q11").Value = 5.96980157435
$ws.Range("R11").Value = 53.72821416915
$ws.Range("S11").Value = 0.01469494580999418
$ws.Range("T11").Value = 0.0189049455143455

# Row 12
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.029143
$ws.Range("H12").Value = 2.058286
$ws.Range("I12").Value = 0.6059291555871485
$ws.Range("J12").Value = 0.5061914997221006
$ws.Range("M12").Value = 55.783591
$ws.Range("N12").Value = 167.350773
$ws.Range("O12").Value = 0.2332214199005771
$ws.Range("P12").Value = 0.2394371967339281
$ws.Range("Q12").Value = 57.409292192513
$ws.Range("R12").Value = 344.455753155078
$ws.Range("S12").Value = 0.1413156580251925
$ws.Range("T12").Value = 0.1212010737040027

# Row 13
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.029143
$ws.Range("H13").Value = 2.058286
$ws.Range("I13").Value = 0.6059291555871485
$ws.Range("J13").Value = 0.5061914997221006
$ws.Range("O13").Value = 0.5297503589663128
$ws.Range("P13").Value = 0.5438691736537713
$ws.Range("Q13").Value = 130.4022296063143
$ws.Range("R13").Value = 782.4133776378859
$ws.Range("S13").Value = 0.3209911876804468
$ws.Range("T13").Value = 0.2753019526644221

# Row 14
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.029143
$ws.Range("H14").Value = 2.058286
$ws.Range("I14").Value = 0.6059291555871485
$ws.Range("J14").Value = 0.5061914997221006
$ws.Range("M14").Value = 23.03749833333333
$ws.Range("N14").Value = 69.112495
$ws.Range("O14").Value = 0.09631574403765399
$ws.Range("P14").Value = 0.09888273454277752
$ws.Range("Q14").Value = 23.70888014726166
$ws.Range("R14").Value = 142.25328088357
$ws.Range("S14").Value = 0.05836051745448362
$ws.Range("T14").Value = 0.05005359969483092

# Row 15
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.029143
$ws.Range("H15").Value = 2.058286
$ws.Range("I15").Value = 0.6059291555871485
$ws.Range("J15").Value = 0.5061914997221006
$ws.Range("M15").Value = 18.627865
$ws.Range("N15").Value = 37.25573
$ws.Range("O15").Value = 0.07787983970082285
$ws.Range("P15").Value = 0.05330365312071852
$ws.Range("Q15").Value = 19.170736869695
$ws.Range("R15").Value = 76.68294747878
$ws.Range("S15").Value = 0.04718966550718208
$ws.Range("T15").Value = 0.02698185611384314

# Row 16
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.029143
$ws.Range("H16").Value = 2.058286
$ws.Range("I16").Value = 0.6059291555871485
$ws.Range("J16").Value = 0.5061914997221006
$ws.Range("M16").Value = 15.028766
$ws.Range("N16").Value = 45.086298
$ws.Range("O16").Value = 0.06283263739463307
$ws.Range("P16").Value = 0.06450724194880479
$ws.Range("Q16").Value = 15.466749327538
$ws.Range("R16").Value = 92.800495965228
$ws.Range("S16").Value = 0.03807212691984351
$ws.Range("T16").Value = 0.0326530175450019
